# Update the "Metrics" sheet figures (B2:B13) with refreshed values.
# Everything downstream (the "today" sheet's B11:B22/E11:E22/F11:F22 formulas,
# which reference Metrics!B2:B13) recalculates automatically.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 173171.38000000003
$metrics.Range("B3").Value  = 152099.04
$metrics.Range("B4").Value  = 53834.33
$metrics.Range("B5").Value  = 7158
$metrics.Range("B6").Value  = 4969417.1300000008
$metrics.Range("B7").Value  = 4194175.7200000007
$metrics.Range("B8").Value  = 1460794.16
$metrics.Range("B9").Value  = 193365
$metrics.Range("B10").Value = 33434798.120000008
$metrics.Range("B11").Value = 31469450.879999999
$metrics.Range("B12").Value = 11742516.200000001
$metrics.Range("B13").Value = 1290995

# Restore the cursor/selection position recorded on the Metrics sheet.
$metrics.Activate()
$metrics.Range("D17").Select()

# Restore the cursor/selection position recorded on the "today" sheet.
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("F4").Select()
